# Apply "custom accuracy" edit:
#  - Row 5 values are rounded to 2 decimal places (matching the target export precision)
#  - Row 6 (the extra data row) is removed entirely, shrinking the used range to A1:AH5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace the 3-decimal raw values with their 2-decimal "custom accuracy" values
$ws.Range("B5").Value  = 14.77
$ws.Range("C5").Value  = 10.98
$ws.Range("D5").Value  = 1.04
$ws.Range("E5").Value  = 32.38
$ws.Range("F5").Value  = 26.14
$ws.Range("G5").Value  = 11.57
$ws.Range("H5").Value  = 41.87
$ws.Range("I5").Value  = 18.03
$ws.Range("J5").Value  = 7.91
$ws.Range("K5").Value  = 11.63
$ws.Range("L5").Value  = 12.98
$ws.Range("M5").Value  = 13.86
$ws.Range("N5").Value  = 3.64
$ws.Range("O5").Value  = 11.65
$ws.Range("P5").Value  = 16.48
$ws.Range("Q5").Value  = 9.949999999999999
$ws.Range("R5").Value  = 0.73
$ws.Range("S5").Value  = 0.66
$ws.Range("T5").Value  = 170.05
$ws.Range("U5").Value  = 32.44
$ws.Range("V5").Value  = 10.76
$ws.Range("W5").Value  = 21.66
$ws.Range("X5").Value  = 11.4
$ws.Range("Y5").Value  = 1.78
$ws.Range("Z5").Value  = 20.72
$ws.Range("AA5").Value = 9.5
$ws.Range("AB5").Value = 8.48
$ws.Range("AC5").Value = 9.960000000000001
$ws.Range("AD5").Value = 13.65
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 37.71
$ws.Range("AG5").Value = 5.99
$ws.Range("AH5").Value = 13.45

# Remove row 6 entirely (data trimmed), which also shrinks the sheet dimension to A1:AH5
$ws.Rows(6).Delete()
